$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for data rows 2-13 from 2023-09-16 (45185) to 2023-10-05 (45204)
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
